# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns for the b4a3cc72... row on the "zh-cn" sheet and the "de-de" sheet
# to reflect the freshly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-21 02:39:21"
$zhcn.Range("H3").Value = "2016-03-21 02:39:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-21 02:39:25"
$dede.Range("H3").Value = "2016-03-21 02:39:47"
